$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.096.17"
$ws.Range("D3").Value = "1.666.40"
$ws.Range("E3").Value = "  -1.32%  "
$ws.Range("D5").Value = "209.49"
$ws.Range("E5").Value = "  -3.60%  "
$ws.Range("E6").Value = "  -1.79%  "
$ws.Range("E7").Value = "  -0.48%  "
$ws.Range("D8").Value = "0.2624"
$ws.Range("E8").Value = "  -3.31%  "
$ws.Range("D9").Value = "0.06288"
$ws.Range("E9").Value = "  -1.72%  "
$ws.Range("D11").Value = "0.07525"
$ws.Range("E11").Value = "  -1.84%  "
$ws.Range("D12").Value = "1.669.52"
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("E13").Value = "  -1.78%  "
$ws.Range("D14").Value = "0.5522"
$ws.Range("E14").Value = "  -4.19%  "
$ws.Range("D15").Value = "66.45"
$ws.Range("E15").Value = "  +0.02%  "
$ws.Range("D16").Value = "0.000007926"
$ws.Range("E16").Value = "  -4.71%  "
$ws.Range("D17").Value = "26.130.94"
$ws.Range("D18").Value = "1.003"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").Value = "4.724"
$ws.Range("E19").Value = "  -3.12%  "
$ws.Range("D20").Value = "186.48"
$ws.Range("E20").Value = "  -1.94%  "
$ws.Range("E21").Value = "  -4.90%  "
$ws.Range("D22").Value = "6.162"
$ws.Range("E22").Value = "  -1.17%  "
$ws.Range("E23").Value = "  -0.50%  "
$ws.Range("D24").Value = "149.63"
$ws.Range("E24").Value = "  +0.79%  "
$ws.Range("D25").Value = "0.1247"
$ws.Range("E25").Value = "  -2.76%  "
$ws.Range("E26").Value = "  -4.47%  "
$ws.Range("D27").Value = "15.89"
$ws.Range("E27").Value = "  +0.64%  "
$ws.Range("D28").Value = "0.06324"
$ws.Range("E28").Value = "  +2.63%  "
$ws.Range("E29").Value = "  -1.25%  "
$ws.Range("E30").Value = "  -3.83%  "
$ws.Range("D31").Value = "3.488"
$ws.Range("E31").Value = "  -2.71%  "
$ws.Range("D32").Value = "3.407"
$ws.Range("E32").Value = "  -4.69%  "
$ws.Range("D33").Value = "1.635"
$ws.Range("E33").Value = "  -2.87%  "
$ws.Range("D34").Value = "0.9966"
$ws.Range("E34").Value = "  -2.95%  "
$ws.Range("D35").Value = "0.6024"
$ws.Range("E35").Value = "  -2.67%  "
$ws.Range("D36").Value = "2.408"
$ws.Range("E36").Value = "  -0.71%  "
$ws.Range("D37").Value = "2.729"
$ws.Range("E37").Value = "  -1.13%  "
$ws.Range("D38").Value = "1.107.65"
$ws.Range("E38").Value = "  +0.23%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.01614"
$ws.Range("E39").Value = "  -1.64%  "
$ws.Range("B40").Value = "FraxShare"
$ws.Range("C40").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D40").Value = "6.080"
$ws.Range("E40").Value = "  -0.90%  "
$ws.Range("D41").Value = "0.8690"
$ws.Range("E41").Value = "  -1.02%  "
$ws.Range("E42").Value = "  -0.93%  "
$ws.Range("D43").Value = "99.79"
$ws.Range("E43").Value = "  -0.96%  "
$ws.Range("D44").Value = "1.818.28"
$ws.Range("E44").Value = "  -1.15%  "
$ws.Range("D45").Value = "0.00000000108"
$ws.Range("E45").Value = "  -2.89%  "
$ws.Range("D46").Value = "55.32"
$ws.Range("E46").Value = "  -3.89%  "
$ws.Range("D47").Value = "1.000"
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("D48").Value = "8.003"
$ws.Range("E48").Value = "  -1.43%  "
$ws.Range("E49").Value = "  -0.87%  "
$ws.Range("D50").Value = "0.4244"
$ws.Range("E50").Value = "  -1.17%  "
$ws.Range("D51").Value = "5.942"
$ws.Range("E51").Value = "  -1.68%  "
